$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update 想去人数 (interest count) values in column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 4458
$wsExhibition.Range("F6").Value = 154
$wsExhibition.Range("F10").Value = 619
$wsExhibition.Range("F12").Value = 192
$wsExhibition.Range("F13").Value = 1240
$wsExhibition.Range("F15").Value = 2867
$wsExhibition.Range("F17").Value = 559

# Sheet "全部类型" (All Types) - same underlying rows, shifted by one row
# because this sheet also includes a "演出" (performance) entry at row 12
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4458
$wsAll.Range("F6").Value = 154
$wsAll.Range("F10").Value = 619
$wsAll.Range("F13").Value = 192
$wsAll.Range("F14").Value = 1240
$wsAll.Range("F16").Value = 2867
$wsAll.Range("F18").Value = 559
